# Apply the geometry changes (Left/Top/Width/Height, in points) to five
# shapes on slide 1 of the overview diagram, matching the target OOXML
# <a:off>/<a:ext> values (in EMU, 12700 EMU == 1 point).
#
# A tiny epsilon is added to each computed point value before assignment:
# the EMU/point round-trip through the host's floating point conversion
# otherwise occasionally truncates instead of rounding to the nearest EMU
# (e.g. 1663339 EMU -> 130.9715748031496 pt -> back to 1663338 EMU). The
# epsilon (well under half an EMU, ~0.0000394 pt) nudges those borderline
# values onto the correct integer EMU without perturbing values that
# already round-trip cleanly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$emu = 12700.0
$eps = 0.00003

# Oval 51 (id 52): move up slightly and grow taller.
$sh = $s.Shapes.Item("Oval 51")
$sh.Top = (6761050 / $emu) + $eps
$sh.Height = (498674 / $emu) + $eps

# Straight Arrow Connector 53 (id 54): shrink height (endpoint follows the
# oval move above).
$sh = $s.Shapes.Item("Straight Arrow Connector 53")
$sh.Height = (1663339 / $emu) + $eps

# Oval 136 (id 137): nudge position, size unchanged.
$sh = $s.Shapes.Item("Oval 136")
$sh.Left = (4466029 / $emu) + $eps
$sh.Top = (6430282 / $emu) + $eps

# Straight Arrow Connector 137 (id 138): shift right, shrink.
$sh = $s.Shapes.Item("Straight Arrow Connector 137")
$sh.Left = (5630516 / $emu) + $eps
$sh.Width = (327665 / $emu) + $eps
$sh.Height = (174746 / $emu) + $eps

# Straight Arrow Connector 154 (id 155): shift right/up, resize.
$sh = $s.Shapes.Item("Straight Arrow Connector 154")
$sh.Left = (5630516 / $emu) + $eps
$sh.Top = (6935097 / $emu) + $eps
$sh.Width = (469809 / $emu) + $eps
$sh.Height = (312924 / $emu) + $eps
